# Aggiornamento fino a 27/05: append new daily rows (256-269) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) of the last existing data row (255) down into the
# new rows so the new date cells in column A keep the same style (s="2") as
# the rest of the date column, and B/C/D keep the unstyled numeric format.
$ws.Range("A255:D255").Copy()
$ws.Range("A256:D269").PasteSpecial(-4122)

# New daily records: date serial, nuovi pos., somma mobile 7gg., somma mobile
# 7gg. per 100mila abitanti
$data = @(
    @(44330, 3,  9,  52.52407353370295),
    @(44331, 4, 11,  64.19608987452582),
    @(44332, 1, 10,  58.36008170411438),
    @(44333, 0, 10,  58.36008170411438),
    @(44334, 1, 10,  58.36008170411438),
    @(44335, 0, 10,  58.36008170411438),
    @(44336, 3, 12,  70.03209804493727),
    @(44337, 0,  9,  52.52407353370295),
    @(44338, 0,  5,  29.18004085205719),
    @(44339, 1,  5,  29.18004085205719),
    @(44340, 0,  5,  29.18004085205719),
    @(44341, 0,  4,  23.34403268164575),
    @(44342, 4,  8,  46.68806536329151),
    @(44343, 0,  5,  29.18004085205719)
)

$r = 256
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
